$d = $word.ActiveDocument

# The last paragraph in the body holds the _GoBack bookmark. Locate it
# dynamically (rather than hard-coding an index) so the script is robust.
$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Insert a brand-new, empty paragraph immediately before the bookmark
# paragraph - this becomes the "Iyanu github" paragraph.
$insertPoint = $bookmarkPara.Range.Duplicate
$insertPoint.Collapse(1)
$insertPoint.InsertParagraphBefore()

# Re-fetch paragraphs after the structural edit; the new paragraph now sits
# right before the (shifted) bookmark paragraph.
$bookmarkIndex = $d.Paragraphs.Count
$newParaIndex = $bookmarkIndex - 1
$newPara = $d.Paragraphs.Item($newParaIndex)

# Build the paragraph text as three separate pieces - "Iyanu", a single
# space, and "github" - mirroring the word-by-word runs (with spell-check
# markers) that Word produced when this text was originally typed.
$wordRange = $newPara.Range.Duplicate
$wordRange.Collapse(1)
$wordRange.InsertAfter("Iyanu")

$wordRange.Collapse(0)
$wordRange.InsertAfter(" ")

$wordRange.Collapse(0)
$wordRange.InsertAfter("github")

# Add a new, empty trailing paragraph right after the bookmark paragraph
# (i.e. just before the final section break).
$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endPoint = $bookmarkPara.Range.Duplicate
$endPoint.Collapse(0)
$endPoint.InsertParagraphAfter()
